# Generate Report for Handback
#
# A handback (df6a600a-6246-4366-bdc7-3466dd5b0682) was produced for both the
# zh-cn and de-de locales, but the handback file turned out to be stale (not
# built from the latest source commit). The report now records:
#   - a link to the (latest) handback markdown file          (col I)
#   - the xliff file that was handed back                    (col J)
#   - the datetime the handback was detected                 (col K)
#   - an error explaining the handback is out of date         (col P)
# The "Error Detail" column is also widened so the message is readable.

$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d2dcab8cfc466730153cd52eb3004a4eba5638b/e2e/df6a600a-6246-4366-bdc7-3466dd5b0682.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a980daef40b2f38ea1a35de72c0329219b2c258b/e2e/df6a600a-6246-4366-bdc7-3466dd5b0682.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d2dcab8cfc466730153cd52eb3004a4eba5638b/e2e/df6a600a-6246-4366-bdc7-3466dd5b0682.md."

# Excel's ColumnWidth property pads the raw stored OOXML width by ~5px worth
# of characters (5/6 of a character for the default Calibri 11 font). The
# other columns in this sheet store a raw width of exactly 40, so back the
# padding out here to land on a stored width of exactly 40 too.
$targetColumnWidth = 40 - (5 / 6)

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I6"), $latestUrl, "", "", "df6a600a-6246-4366-bdc7-3466dd5b0682.md")
$wsZh.Range("J6").Value = $wsZh.Range("G6").Value2
$wsZh.Range("K6").Value = "2016-08-20 12:46:40"
$wsZh.Range("P6").Value = $errorDetail

$wsZh.Columns.Item(16).ColumnWidth = $targetColumnWidth

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I6"), $latestUrl, "", "", "df6a600a-6246-4366-bdc7-3466dd5b0682.md")
$wsDe.Range("J6").Value = $wsDe.Range("G6").Value2
$wsDe.Range("K6").Value = "2016-08-20 12:46:46"
$wsDe.Range("P6").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = $targetColumnWidth
